# Test case successful for unicode and special characters
#
# Adds a new "TU05" test-case row (row 6) to the test-cases sheet, styled
# to match the existing orange/"pass" row above it (row 5), re-colors
# that accent fill to a salmon/pink tone, clears the old fill from row 5,
# grows row 6's height and moves the active selection down to B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New test-case values for row 6 (appended to the shared string table
#    in left-to-right order, matching how Excel would record new text).
$ws.Range("A6").Value = "TU05"
$ws.Range("B6").Value = "Check response when passing a city name with unicode in it."
$ws.Range("C6").Value = "./weather São Paulo"
$ws.Range("D6").Value = "Enter a city with special character like São Paulo"
$ws.Range("E6").Value = "Converted to plain alphabet and then pushed to array"
$ws.Range("F6").Value = "Pass"

# 2) Row 6 grows to fit (matches the taller "pass" rows above it).
$ws.Rows(6).RowHeight = 30

# 3) Give row 6 the same visual treatment as the other data rows: thin
#    border all around, centered alignment, and the accent fill/font
#    combo used for "Pass" rows - but recolored to the new salmon tone.
foreach ($col in @("A", "B", "C", "D", "E", "F")) {
    $cell = $ws.Range($col + "6")

    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Color = 0x000000

    $cell.Interior.Color = 0x9AA1F7
    $cell.Interior.PatternColor = 0x8080FF

    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2

    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Columns B, D and E hold longer descriptive text, so they wrap like the
# corresponding columns in the rows above.
foreach ($col in @("B", "D", "E")) {
    $ws.Range($col + "6").WrapText = $true
}

# 4) The accent fill used on row 5 is recolored (from orange to salmon)
#    and handed off to row 6; row 5 goes back to no fill (preserving the
#    wrap setting on columns B and D, which the fill change would
#    otherwise reset).
foreach ($col in @("A", "B", "C", "D", "E", "F")) {
    $ws.Range($col + "5").Interior.Pattern = -4142
}
foreach ($col in @("B", "D")) {
    $ws.Range($col + "5").WrapText = $true
}

# 5) Move the active selection to where editing ended up (B17).
$ws.Range("B17").Select()
